$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B (so Time column is inserted after Date)
$ws.Columns("B").Insert()

# Set header for new Time column
$ws.Range("B1").Value = "Time"

# Set Time values for rows 2 and 3
$ws.Range("B2").Value = "04:48"
$ws.Range("B3").Value = "04:48"

# Update Date format in column A for rows 2 and 3
$ws.Range("A2").Value = "19-01-26"
$ws.Range("A3").Value = "19-01-26"
